$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "Find duplicates in O(n) time and O(1) extra space"
$ws.Range("B16").Value = "FindDuplicate"

$ws.Range("B16").Select()
